# Update Active_Outages.xlsx - 6/19/2025, 1:56:17 PM
#
# 1) Refresh the "Elapsed Duration(Hrs)" column (G) on several outage rows
#    across the region sheets (R1, R2, R4, R5, R6) with newer elapsed times.
# 2) Append a new outage row (row 6) to sheet "R1", duplicating the
#    existing JED0125 / Generator-SG / "Good+In progress" outage row that
#    already exists on sheet "R2".

$wb = $excel.ActiveWorkbook

# --- 1) Elapsed Duration(Hrs) refresh -------------------------------------
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3951:10:15"
$ws1.Range("G3").Value = "90:42:53"
$ws1.Range("G4").Value = "113:42:53"

$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12132:33:56"
$ws2.Range("G3").Value = "3262:17:25"
$ws2.Range("G4").Value = "500:28:59"

$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2978:23:45"
$ws4.Range("G3").Value = "205:36:00"
$ws4.Range("G4").Value = "93:48:25"
$ws4.Range("G5").Value = "91:25:58"

$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "452:22:44"

$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "92:55:02"

# --- 2) Append new outage row to R1 (copy of R2 row 6) --------------------
$ws1.Range("A6").Value = ""
$ws1.Range("B6").Value = "R4"
$ws1.Range("C6").Value = ""
$ws1.Range("D6").Value = "JED0125"
$ws1.Range("E6").Value = ""
$ws1.Range("F6").Value = ""
$ws1.Range("G6").Value = ""
$ws1.Range("H6").Value = ""
$ws1.Range("I6").Value = "Generator-SG"
$ws1.Range("J6").Value = "Good+In progress"
$ws1.Range("K6").Value = ""
$ws1.Range("L6").Value = "Latis"
